$wb = $excel.ActiveWorkbook

# --- RUNMANAGER sheet: set C3:C6 to "yes", update selection ---
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws1.Range("C3").Value = "yes"
$ws1.Range("C4").Value = "yes"
$ws1.Range("C5").Value = "yes"
$ws1.Range("C6").Value = "yes"
$ws1.Range("C2:C6").Select()

# --- DATA sheet: set B3:B6 to "yes", update selection, and make it the active sheet ---
$ws2 = $wb.Worksheets.Item("DATA")
$ws2.Range("B3").Value = "yes"
$ws2.Range("B4").Value = "yes"
$ws2.Range("B5").Value = "yes"
$ws2.Range("B6").Value = "yes"
$ws2.Activate()
$ws2.Range("B2:B6").Select()
